# Fixing GAB Template to not need background slide
# Adds a full-bleed blue rectangle to the slide so the template no longer
# relies on an external background slide, and marks the headline text run
# as "clean" (dirty="0") the way PowerPoint does after an edit pass.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Register an (empty) slide-guide list extension on the presentation,
#        mirroring what PowerPoint writes once guides have been touched.
try {
    $null = $p.Guides.Add(1, 360)
} catch {
}

# --- 2. Insert the new background rectangle ("Rectangle 1") behind every
#        other shape already on the slide.
#        Position/size below are taken straight from the target EMU values,
#        converted to points (1 pt = 12700 EMU):
#          off  x=0        y=857250   -> 0, 67.5
#          ext  cx=9144000 cy=5143500 -> 720, 405
$rect = $s.Shapes.AddShape(1, 0, 67.5, 720, 405)
$rect.Name = "Rectangle 1"

# Push it to the very back of the z-order so it acts as a background.
$rect.ZOrder(1)

# Solid fill 0064A2, no outline.
$rect.Fill.ForeColor.RGB = 0xA26400
$rect.Line.Visible = $false

# Centered text body / paragraph alignment, matching the target markup.
$rect.TextFrame.VerticalAnchor = 3
$rect.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- 3. Touch up the headline run on "Shape 77" so it is marked dirty="0".
$headline = $s.Shapes.Item("Shape 77")
$headline.TextFrame.TextRange.Text = "LIVING HIS MESSAGE - Prayers for Global Peace"
